$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Wins / Losses / Ties, styled like the rest of row 1 (copy AC1's
# formatting - bold font + border + alignment - onto the new header cells)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record columns for every data row (2-46): Wins=87, Losses=75, Ties=0
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 87
    $ws.Cells.Item($r, 31).Value = 75
    $ws.Cells.Item($r, 32).Value = 0
}
